$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) New header cells T1 (fuente_tiempos) and U1 (estado_datos), matching header style ---
$ws.Range("S1").Copy() | Out-Null
$ws.Range("T1:U1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false
$ws.Range("T1").Value = "fuente_tiempos"
$ws.Range("U1").Value = "estado_datos"

# --- 2) Existing rows 2-14: add T (fuente_tiempos="score") and U (estado_datos="OK") ---
for ($r = 2; $r -le 14; $r++) {
    $ws.Range("T$r").Value = "score"
    $ws.Range("U$r").Value = "OK"
}

# --- 3) New fixture rows 15-20 ---

# Row 15
$ws.Range("A15").NumberFormat = "@"
$ws.Range("A15").Value = "2025-08-09"
$ws.Range("A15").ClearFormats()
$ws.Range("B15").Value = "FC ST. Gallen"
$ws.Range("C15").Value = "FC Winterthur"
$ws.Range("D15").Value = 5
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1382308
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = 5
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 1
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = 3
$ws.Range("N15").Value = 0
$ws.Range("O15").Value = 2
$ws.Range("P15").Value = 0
$ws.Range("Q15").Value = 36
$ws.Range("R15").Value = 64
$ws.Range("S15").Value = "L"
$ws.Range("T15").Value = "score"
$ws.Range("U15").Value = "OK"

# Row 16
$ws.Range("A16").NumberFormat = "@"
$ws.Range("A16").Value = "2025-08-09"
$ws.Range("A16").ClearFormats()
$ws.Range("B16").Value = "FC Luzern"
$ws.Range("C16").Value = "FC Thun"
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 1382309
$ws.Range("G16").Value = 6
$ws.Range("H16").Value = 5
$ws.Range("I16").Value = 2
$ws.Range("J16").Value = 2
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = 1
$ws.Range("N16").Value = 1
$ws.Range("O16").Value = 0
$ws.Range("P16").Value = 1
$ws.Range("Q16").Value = 55
$ws.Range("R16").Value = 45
$ws.Range("S16").Value = "V"
$ws.Range("T16").Value = "score"
$ws.Range("U16").Value = "OK"

# Row 17
$ws.Range("A17").NumberFormat = "@"
$ws.Range("A17").Value = "2025-08-10"
$ws.Range("A17").ClearFormats()
$ws.Range("B17").Value = "BSC Young Boys"
$ws.Range("C17").Value = "FC Sion"
$ws.Range("D17").Value = 0
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 1382310
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 5
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 3
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 1
$ws.Range("M17").Value = 0
$ws.Range("N17").Value = 0
$ws.Range("O17").Value = 0
$ws.Range("P17").Value = 0
$ws.Range("Q17").Value = 58
$ws.Range("R17").Value = 42
$ws.Range("S17").Value = "E"
$ws.Range("T17").Value = "score"
$ws.Range("U17").Value = "OK"

# Row 18
$ws.Range("A18").NumberFormat = "@"
$ws.Range("A18").Value = "2025-08-10"
$ws.Range("A18").ClearFormats()
$ws.Range("B18").Value = "FC Lugano"
$ws.Range("C18").Value = "FC Basel 1893"
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = 1382312
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 6
$ws.Range("I18").Value = 2
$ws.Range("J18").Value = 3
$ws.Range("K18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = 2
$ws.Range("N18").Value = 0
$ws.Range("O18").Value = 1
$ws.Range("P18").Value = 1
$ws.Range("Q18").Value = 50
$ws.Range("R18").Value = 50
$ws.Range("S18").Value = "L"
$ws.Range("T18").Value = "score"
$ws.Range("U18").Value = "OK"

# Row 19
$ws.Range("A19").NumberFormat = "@"
$ws.Range("A19").Value = "2025-08-10"
$ws.Range("A19").ClearFormats()
$ws.Range("B19").Value = "Lausanne"
$ws.Range("C19").Value = "FC Zurich"
$ws.Range("D19").Value = 1
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = 1382311
$ws.Range("G19").Value = 8
$ws.Range("H19").Value = 6
$ws.Range("I19").Value = 3
$ws.Range("J19").Value = 3
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 0
$ws.Range("M19").Value = 1
$ws.Range("N19").Value = 1
$ws.Range("O19").Value = 0
$ws.Range("P19").Value = 1
$ws.Range("Q19").Value = 52
$ws.Range("R19").Value = 48
$ws.Range("S19").Value = "V"
$ws.Range("T19").Value = "score"
$ws.Range("U19").Value = "OK"

# Row 20
$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "2025-08-10"
$ws.Range("A20").ClearFormats()
$ws.Range("B20").Value = "Servette FC"
$ws.Range("C20").Value = "Grasshoppers"
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 1
$ws.Range("F20").Value = 1382307
$ws.Range("G20").Value = 8
$ws.Range("H20").Value = 5
$ws.Range("I20").Value = 2
$ws.Range("J20").Value = 1
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("M20").Value = 0
$ws.Range("N20").Value = 1
$ws.Range("O20").Value = 1
$ws.Range("P20").Value = 0
$ws.Range("Q20").Value = 59
$ws.Range("R20").Value = 41
$ws.Range("S20").Value = "E"
$ws.Range("T20").Value = "score"
$ws.Range("U20").Value = "OK"

Write-Output "Edit complete"
